# Implement hh-based ACF around TB cases detected through PCF
#
# Changes:
#  1. "n_hh" (B1 on the "constant" sheet) parameter bumped from 1000 to 10000
#  2. A new parameter row is appended to the "constant" sheet:
#       A79 = "hh_based_acf_coverage_perc"
#       B79 = 0
#       C79 = "float"
#     formatted the same way as the other rows in that block (copy format
#     from the row above so fill/border style "1" carries over).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constant")

# 1. Update the existing parameter value
$ws.Range("B1").Value = 10000

# 2. Append the new "hh_based_acf_coverage_perc" parameter row, reusing the
#    formatting of the row directly above it (row 44) before writing values.
$ws.Range("A44:C44").Copy()
$ws.Range("A79:C79").PasteSpecial(-4122)

$ws.Range("A79").Value = "hh_based_acf_coverage_perc"
$ws.Range("B79").Value = 0
$ws.Range("C79").Value = "float"

# Leave the new cell selected, matching the saved view state.
$ws.Range("B79").Select()
